# Auto-generated script to update Sheets via scheduled runner
# Applies computed market price / profit refresh values to the Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value2 = 475.875
$ws.Range("I38").Value2 = 200.625
$ws.Range("J38").Value2 = 751.125
$ws.Range("K38").Value2 = 601.875
$ws.Range("L38").Value2 = 2253.375
$ws.Range("M38").Value2 = -229.875
$ws.Range("N38").Value2 = -2997.375
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 8173.3335
$ws.Range("K51").Value2 = 0
$ws.Range("L51").Value2 = 8173.3335
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value2 = -9141.333500000001
$ws.Range("H113").Value2 = 3092.9792
$ws.Range("J113").Value2 = 4866.278
$ws.Range("L113").Value2 = 4866.278
$ws.Range("N113").Value2 = -11374.278
$ws.Range("H129").Value2 = 804.0808
$ws.Range("I129").Value2 = 0
$ws.Range("J129").Value2 = 804.0808
$ws.Range("K129").Value2 = 0
$ws.Range("L129").Value2 = 2412.2424
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value2 = -12412.2424
$ws.Range("H132").Value2 = 2819.7097
$ws.Range("I132").Value2 = 2996.6072
$ws.Range("J132").Value2 = 1168.6666
$ws.Range("K132").Value2 = 8989.821599999999
$ws.Range("L132").Value2 = 3505.9998
$ws.Range("M132").Value2 = -6459.821599999999
$ws.Range("N132").Value2 = -8565.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1377.0883
$ws.Range("I2").Value2 = 858.65216
$ws.Range("J2").Value2 = 2461.0908
$ws.Range("K2").Value2 = 858.65216
$ws.Range("L2").Value2 = 2461.0908
$ws.Range("M2").Value2 = -745.65216
$ws.Range("N2").Value2 = -2687.0908
$ws.Range("H32").Value2 = 3527.075
$ws.Range("I32").Value2 = 2938.4517
$ws.Range("J32").Value2 = 5554.5557
$ws.Range("K32").Value2 = 2938.4517
$ws.Range("L32").Value2 = 5554.5557
$ws.Range("M32").Value2 = -2651.4517
$ws.Range("N32").Value2 = -6128.5557
$ws.Range("H110").Value2 = 1826.25
$ws.Range("J110").Value2 = 1419.8334
$ws.Range("L110").Value2 = 1419.8334
$ws.Range("N110").Value2 = -5509.8334
$ws.Range("H116").Value2 = 1377.0883
$ws.Range("I116").Value2 = 858.65216
$ws.Range("J116").Value2 = 2461.0908
$ws.Range("K116").Value2 = 858.65216
$ws.Range("L116").Value2 = 2461.0908
$ws.Range("M116").Value2 = 1435.34784
$ws.Range("N116").Value2 = -7049.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1377.0883
$ws.Range("I3").Value2 = 858.65216
$ws.Range("J3").Value2 = 2461.0908
$ws.Range("K3").Value2 = 858.65216
$ws.Range("L3").Value2 = 2461.0908
$ws.Range("M3").Value2 = -744.65216
$ws.Range("N3").Value2 = -2689.0908
$ws.Range("H94").Value2 = 4162.7393
$ws.Range("I94").Value2 = 2171
$ws.Range("K94").Value2 = 2171
$ws.Range("M94").Value2 = -1720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value2 = 0
$ws.Range("J11").Value2 = 0
$ws.Range("L11").Value2 = 0
$ws.Range("N11").ClearContents()
$ws.Range("H31").Value2 = 10796
$ws.Range("I31").Value2 = 12149.417
$ws.Range("J31").Value2 = 3835.5715
$ws.Range("K31").Value2 = 12149.417
$ws.Range("L31").Value2 = 3835.5715
$ws.Range("M31").Value2 = -11854.417
$ws.Range("N31").Value2 = -4425.5715
$ws.Range("H34").Value2 = 10796
$ws.Range("I34").Value2 = 12149.417
$ws.Range("J34").Value2 = 3835.5715
$ws.Range("K34").Value2 = 12149.417
$ws.Range("L34").Value2 = 3835.5715
$ws.Range("M34").Value2 = -11947.417
$ws.Range("N34").Value2 = -4239.5715
$ws.Range("H99").Value2 = 4911.3
$ws.Range("I99").Value2 = 3622.2856
$ws.Range("J99").Value2 = 7919
$ws.Range("K99").Value2 = 3622.2856
$ws.Range("L99").Value2 = 7919
$ws.Range("M99").Value2 = -2124.2856
$ws.Range("N99").Value2 = -10915
$ws.Range("H126").Value2 = 4911.3
$ws.Range("I126").Value2 = 3622.2856
$ws.Range("J126").Value2 = 7919
$ws.Range("K126").Value2 = 10866.8568
$ws.Range("L126").Value2 = 23757
$ws.Range("M126").Value2 = -8396.856800000001
$ws.Range("N126").Value2 = -28697
$ws.Range("H132").Value2 = 20726.678
$ws.Range("I132").Value2 = 22721.75
$ws.Range("K132").Value2 = 68165.25
$ws.Range("M132").Value2 = -65635.25
$ws.Range("H134").Value2 = 1096.6666
$ws.Range("I134").Value2 = 826.1053000000001
$ws.Range("K134").Value2 = 2478.3159
$ws.Range("M134").Value2 = 56.68409999999994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value2 = 93.666664
$ws.Range("I6").Value2 = 74.72727
$ws.Range("J6").Value2 = 302
$ws.Range("K6").Value2 = 224.18181
$ws.Range("L6").Value2 = 906
$ws.Range("M6").Value2 = -111.18181
$ws.Range("N6").Value2 = -1132
$ws.Range("H33").Value2 = 266.33334
$ws.Range("J33").Value2 = 300
$ws.Range("L33").Value2 = 1800
$ws.Range("N33").Value2 = -2366
$ws.Range("H113").Value2 = 20505.8
$ws.Range("I113").Value2 = 33709.668
$ws.Range("J113").Value2 = 700
$ws.Range("K113").Value2 = 101129.004
$ws.Range("L113").Value2 = 2100
$ws.Range("M113").Value2 = -98959.00399999999
$ws.Range("N113").Value2 = -6440
$ws.Range("H131").Value2 = 709.63
$ws.Range("I131").Value2 = 307.5
$ws.Range("J131").Value2 = 726.38544
$ws.Range("K131").Value2 = 922.5
$ws.Range("L131").Value2 = 2179.15632
$ws.Range("M131").Value2 = 4117.5
$ws.Range("N131").Value2 = -12259.15632

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value2 = 4269.5
$ws.Range("I53").Value2 = 4026
$ws.Range("K53").Value2 = 4026
$ws.Range("M53").Value2 = -3395
$ws.Range("H113").Value2 = 3788.7778
$ws.Range("I113").Value2 = 3333.3333
$ws.Range("K113").Value2 = 3333.3333
$ws.Range("M113").Value2 = -1163.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 5112.3687
$ws.Range("I61").Value2 = 2255.5557
$ws.Range("K61").Value2 = 2255.5557
$ws.Range("M61").Value2 = -2053.5557
$ws.Range("H113").Value2 = 5112.3687
$ws.Range("I113").Value2 = 2255.5557
$ws.Range("K113").Value2 = 2255.5557
$ws.Range("M113").Value2 = -85.55569999999989
$ws.Range("H132").Value2 = 1893.3182
$ws.Range("I132").Value2 = 870.8333
$ws.Range("J132").Value2 = 3120.3
$ws.Range("K132").Value2 = 2612.4999
$ws.Range("L132").Value2 = 9360.900000000001
$ws.Range("M132").Value2 = -82.4998999999998
$ws.Range("N132").Value2 = -14420.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value2 = 31658.666
$ws.Range("J27").Value2 = 31658.666
$ws.Range("L27").Value2 = 31658.666
$ws.Range("N27").Value2 = -31796.666
$ws.Range("H113").Value2 = 3862604
$ws.Range("I113").Value2 = 2650
$ws.Range("J113").Value2 = 9009209
$ws.Range("K113").Value2 = 7950
$ws.Range("L113").Value2 = 27027627
$ws.Range("M113").Value2 = -5780
$ws.Range("N113").Value2 = -27031967
$ws.Range("H132").Value2 = 3119
$ws.Range("I132").Value2 = 2687.8
$ws.Range("J132").Value2 = 3837.6667
$ws.Range("K132").Value2 = 8063.400000000001
$ws.Range("L132").Value2 = 11513.0001
$ws.Range("M132").Value2 = -5533.400000000001
$ws.Range("N132").Value2 = -16573.0001
$ws.Range("H136").Value2 = 1298.9459
$ws.Range("I136").Value2 = 867.73914
$ws.Range("J136").Value2 = 2007.3572
$ws.Range("K136").Value2 = 2603.21742
$ws.Range("L136").Value2 = 6022.071599999999
$ws.Range("M136").Value2 = -53.21741999999995
$ws.Range("N136").Value2 = -11122.0716
